$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value (kept as Text to preserve original formatting)
$updates = @{
    "D2" = "246.27"
    "E2" = "-0.12%"
    "D3" = "29.74"
    "E3" = "-1.34%"
    "D4" = "5.157"
    "E4" = "-0.62%"
    "D5" = "0.05765"
    "E5" = "0.57%"
    "D6" = "6.661"
    "E6" = "1.00%"
    "D7" = "3.227"
    "E7" = "5.98%"
    "D8" = "0.8502"
    "E8" = "-1.06%"
    "D9" = "0.8561"
    "E9" = "-1.95%"
    "E10" = "2.40%"
    "D11" = "0.07092"
    "E11" = "-0.07%"
    "D12" = "0.03249"
    "E12" = "11.14%"
    "D13" = "0.09373"
    "E13" = "-0.25%"
    "D14" = "0.001539"
    "E14" = "1.47%"
    "D15" = "0.0005966"
    "E15" = "-94.15%"
    "D16" = "0.005888"
    "E16" = "-4.89%"
    "D17" = "3.522"
    "E17" = "0.45%"
    "D18" = "2.213"
    "E18" = "-2.90%"
    "E19" = "-0.51%"
    "D20" = "0.03386"
    "E20" = "2.93%"
    "E21" = "0.31%"
    "D22" = "3.504"
    "E22" = "-2.93%"
    "E23" = "2.16%"
    "E24" = "-0.55%"
    "D25" = "0.001227"
    "E25" = "1.06%"
    "E26" = "-7.87%"
    "E27" = "1.65%"
    "E28" = "4.14%"
    "D40" = "0.03751"
    "E40" = "-0.87%"
    "D41" = "0.1071"
    "E41" = "0.02%"
    "D42" = "0.002198"
    "E42" = "-0.04%"
    "D43" = "0.002947"
    "E43" = "-48.27%"
    "D44" = "0.009937"
    "E44" = "-0.76%"
    "D45" = "0.00005475"
    "E45" = "7.35%"
    "E46" = "-0.04%"
    "D47" = "0.07096"
    "E47" = "-20.23%"
    "D48" = "0.002466"
    "E48" = "-10.87%"
    "E49" = "-0.04%"
    "E50" = "-0.04%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
